# 11.6.1.1a.xlsx update
#  - insert a new row below the title row with a "(in percent)" style
#    sub-header (Kyrgyz / Russian / English) in columns A:C
#  - add two new data columns (L, M) for years 2021 and 2022
#  - move the selection cursor

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row 2 (shifts the old rows 2-9 down to 3-10) and fill the
#    new "(в процентах)" sub-header cells.
# ---------------------------------------------------------------------------
$ws.Rows("2:2").Insert()

$ws.Range("A2").Value = "(пайыз менен)"
$ws.Range("B2").Value = "(в процентах)"
$ws.Range("C2").Value = "(in percent)"

# Build the italic 8pt centered/wrapped look on a pristine scratch cell (so
# we start from the sheet's plain default font instead of inheriting the
# bold title-row font), then stamp that format onto the new header cells.
$scratch = $ws.Range("Z1")
$scratch.Font.Size = 8
$scratch.Font.Italic = $true
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4108
$scratch.WrapText = $true

$scratch.Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
$scratch.Clear()

# ---------------------------------------------------------------------------
# 2) Add the 2021 / 2022 data columns (L, M), reusing each row's existing
#    number format / border (copied from column K).
# ---------------------------------------------------------------------------
$ws.Range("K4").Copy()
$ws.Range("L4:M4").PasteSpecial(-4122)

$ws.Range("K5").Copy()
$ws.Range("L5:M5").PasteSpecial(-4122)

$ws.Range("K6").Copy()
$ws.Range("L6:M6").PasteSpecial(-4122)

$ws.Range("K7").Copy()
$ws.Range("L7:M7").PasteSpecial(-4122)

$ws.Range("K8").Copy()
$ws.Range("L8:M8").PasteSpecial(-4122)

$ws.Range("K9").Copy()
$ws.Range("L9:M9").PasteSpecial(-4122)

$ws.Range("L4").Value = 2021
$ws.Range("M4").Value = 2022

$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0

$ws.Range("L6").Value = 58.405380200320216
$ws.Range("M6").Value = 48.6

$ws.Range("L7").Value = 11.673077354810609
$ws.Range("M7").Value = 20.2

$ws.Range("L8").Value = 22.564920591204277
$ws.Range("M8").Value = 22.9

$ws.Range("L9").Value = 7.3566218536648895
$ws.Range("M9").Value = 8.3000000000000007

# ---------------------------------------------------------------------------
# 3) Move the active selection.
# ---------------------------------------------------------------------------
$ws.Range("N7").Select() | Out-Null
